$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 607, shifting existing rows 607:714 down to 608:715.
$ws.Rows(607).Insert()

# Populate the newly inserted row 607 with the new weekly price record.
$ws.Cells.Item(607, 1).Value  = 3
$ws.Cells.Item(607, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(607, 3).Value  = "Coquimbo"
$ws.Cells.Item(607, 4).Value  = 45180
$ws.Cells.Item(607, 5).Value  = 5
$ws.Cells.Item(607, 6).Value  = 100112017
$ws.Cells.Item(607, 7).Value  = "Apio"
$ws.Cells.Item(607, 8).Value  = "Americana (o)"
$ws.Cells.Item(607, 9).Value  = "Primera"
$ws.Cells.Item(607, 10).Value = 200
$ws.Cells.Item(607, 11).Value = 9000
$ws.Cells.Item(607, 12).Value = 9500
$ws.Cells.Item(607, 13).Value = 9200
$ws.Cells.Item(607, 14).Value = "`$/docena de matas"
$ws.Cells.Item(607, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(607, 16).Value = 1533
$ws.Cells.Item(607, 17).Value = 6
$ws.Cells.Item(607, 18).Value = "Hortaliza"
